$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F1").Value = "sd_Y"

$values = @(
    0.00853036609572474,
    0.0464544794571031,
    0.103670497317918,
    0.0852520525247769,
    0,
    0.175508950274643,
    0.156580542054294,
    0.116960336174683,
    0.0185841285177258,
    0.105162378480623,
    0.0978254310552359,
    0.0471294735149011,
    0.0159026501014877,
    0.204078309925927,
    0.219477830982072,
    0
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $values[$i]
}
